# CS133JS Lab04 Rubric - "New and revised rubrics"
# Update the point values on both the "Rubric" and "Score" sheets,
# and restore the selection/view state shown in the target workbook.

$wb = $excel.ActiveWorkbook
$wsRubric = $wb.Worksheets.Item("Rubric")
$wsScore  = $wb.Worksheets.Item("Score")

# ----- Rubric sheet: revised point values in column D -----
$wsRubric.Range("D9").Value  = 2
$wsRubric.Range("D11").Value = 3
$wsRubric.Range("D12").Value = 3
$wsRubric.Range("D13").Value = 4
$wsRubric.Range("D15").Value = 4
$wsRubric.Range("D16").Value = 5

# ----- Score sheet: same revised point values, columns D and E -----
$wsScore.Range("D9").Value  = 2
$wsScore.Range("E9").Value  = 2

$wsScore.Range("D11").Value = 3
$wsScore.Range("E11").Value = 3

$wsScore.Range("D12").Value = 3
$wsScore.Range("E12").Value = 3

$wsScore.Range("D13").Value = 4
$wsScore.Range("E13").Value = 4

$wsScore.Range("D15").Value = 4
$wsScore.Range("E15").Value = 4

$wsScore.Range("D16").Value = 5
$wsScore.Range("E16").Value = 5

# ----- Restore view/selection state -----
# Score sheet: zoomed to 160%, scrolled so row 4 is at top, G16 selected
$wsScore.Activate()
$scoreWindow = $excel.ActiveWindow
$scoreWindow.Zoom = 160
$scoreWindow.ScrollRow = 4
$scoreWindow.ScrollColumn = 1
$wsScore.Range("G16").Select()

# Rubric sheet stays the active/selected tab, with H21 selected
$wsRubric.Activate()
$wsRubric.Range("H21").Select()
